# Test Page Object Model Commit
# Appends numeric suffixes to the sample FirstName/Username values used by
# the Selenium Page-Object-Model tests (RegisterSheet + mirrored LoginSheet),
# clears the stray fill style left on the last Password cell of each sheet,
# and restores the originally-selected cell on each sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RegisterSheet")
$ws2 = $wb.Worksheets.Item("LoginSheet")

# --- RegisterSheet: FirstName (A) column gets a numeric suffix ---
$ws1.Range("A2").Value = "Sai1"
$ws1.Range("A3").Value = "Prabal2"
$ws1.Range("A4").Value = "Abhas3"
$ws1.Range("A5").Value = "Bhagya4"

# --- RegisterSheet: Username (I) column gets a numeric suffix ---
$ws1.Range("I2").Value = "SaiSundar101"
$ws1.Range("I3").Value = "Ghosh102"
$ws1.Range("I4").Value = "Abhas103"
$ws1.Range("I5").Value = "Bhagya104"

# --- LoginSheet: Username column (A) mirrors RegisterSheet's Username column ---
$ws2.Range("A2").Value = "SaiSundar101"
$ws2.Range("A3").Value = "Ghosh102"
$ws2.Range("A4").Value = "Abhas103"
$ws2.Range("A5").Value = "Bhagya104"

# --- Drop the leftover "applyFill" style from the last Password cell on each sheet ---
$ws1.Range("J5").Style = "Normal"
$ws2.Range("B5").Style = "Normal"

# --- Restore the selected cell recorded on each sheet ---
$ws1.Activate()
$ws1.Range("F8").Select() | Out-Null

$ws2.Activate()
$ws2.Range("F2").Select() | Out-Null
